$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate current row 2 (copies its formatting, incl. the date style on column D)
# into a freshly inserted row 3, pushing all existing data rows (2-10) down to (3-11).
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).Insert()

# Overwrite row 2 with the new weekly record.
$ws.Cells.Item(2, 1).Value = 11
$ws.Cells.Item(2, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(2, 3).Value = "Bíobío"
$ws.Cells.Item(2, 4).Value = 44503
$ws.Cells.Item(2, 5).Value = 8
$ws.Cells.Item(2, 6).Value = 100112022
$ws.Cells.Item(2, 7).Value = "Arveja Verde"
$ws.Cells.Item(2, 8).Value = "Perfection"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 200
$ws.Cells.Item(2, 11).Value = 15000
$ws.Cells.Item(2, 12).Value = 16000
$ws.Cells.Item(2, 13).Value = 15500
$ws.Cells.Item(2, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(2, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(2, 16).Value = 620
$ws.Cells.Item(2, 17).Value = 25
$ws.Cells.Item(2, 18).Value = "Hortaliza"

$wb.Save()
